$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Cilantro" series. Insert a
# fresh row at position 40 (shifting the existing rows 40-94 down to 41-95)
# and populate it with the new observation.
$ws.Rows(40).Insert()

$ws.Cells.Item(40, 1).Value  = 5
$ws.Cells.Item(40, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(40, 3).Value  = "Maule"
$ws.Cells.Item(40, 4).Value  = 45079
$ws.Cells.Item(40, 5).Value  = 7
$ws.Cells.Item(40, 6).Value  = 100112040
$ws.Cells.Item(40, 7).Value  = "Cilantro"
$ws.Cells.Item(40, 8).Value  = "Sin especificar"
$ws.Cells.Item(40, 9).Value  = "Primera"
$ws.Cells.Item(40, 10).Value = 150
$ws.Cells.Item(40, 11).Value = 7000
$ws.Cells.Item(40, 12).Value = 7000
$ws.Cells.Item(40, 13).Value = 7000
$ws.Cells.Item(40, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(40, 15).Value = "Región del Maule"
$ws.Cells.Item(40, 16).Value = 194
$ws.Cells.Item(40, 17).Value = 36
$ws.Cells.Item(40, 18).Value = "Hortaliza"
